$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style (including number format, border, font, alignment) of A255
# into the new date cells so they pick up the same cellXfs entry (s="2"),
# then overwrite the values straight after - Copy(Destination) keeps the xf index
# stable instead of fabricating a brand new style entry.
$ws.Range("A255").Copy($ws.Range("A256"))
$ws.Range("A256").Value = 44330
$ws.Range("B256").Value = 3
$ws.Range("C256").Value = 43
$ws.Range("D256").Value = 239.1812214929358

$ws.Range("A255").Copy($ws.Range("A257"))
$ws.Range("A257").Value = 44331
$ws.Range("B257").Value = 1
$ws.Range("C257").Value = 40
$ws.Range("D257").Value = 222.4941595283124

$ws.Range("A255").Copy($ws.Range("A258"))
$ws.Range("A258").Value = 44332
$ws.Range("B258").Value = 1
$ws.Range("C258").Value = 30
$ws.Range("D258").Value = 166.8706196462343

$ws.Range("A255").Copy($ws.Range("A259"))
$ws.Range("A259").Value = 44333
$ws.Range("B259").Value = 2
$ws.Range("C259").Value = 31
$ws.Range("D259").Value = 172.4329736344421

$ws.Range("A255").Copy($ws.Range("A260"))
$ws.Range("A260").Value = 44334
$ws.Range("B260").Value = 1
$ws.Range("C260").Value = 21
$ws.Range("D260").Value = 116.809433752364

$ws.Range("A255").Copy($ws.Range("A261"))
$ws.Range("A261").Value = 44335
$ws.Range("B261").Value = 1
$ws.Range("C261").Value = 20
$ws.Range("D261").Value = 111.2470797641562

$ws.Range("A255").Copy($ws.Range("A262"))
$ws.Range("A262").Value = 44336
$ws.Range("B262").Value = 4
$ws.Range("C262").Value = 13
$ws.Range("D262").Value = 72.31060184670152

$ws.Range("A255").Copy($ws.Range("A263"))
$ws.Range("A263").Value = 44337
$ws.Range("B263").Value = 0
$ws.Range("C263").Value = 10
$ws.Range("D263").Value = 55.6235398820781

$ws.Range("A255").Copy($ws.Range("A264"))
$ws.Range("A264").Value = 44338
$ws.Range("B264").Value = 1
$ws.Range("C264").Value = 10
$ws.Range("D264").Value = 55.6235398820781

$ws.Range("A255").Copy($ws.Range("A265"))
$ws.Range("A265").Value = 44339
$ws.Range("B265").Value = 0
$ws.Range("C265").Value = 9
$ws.Range("D265").Value = 50.06118589387028

$ws.Range("A255").Copy($ws.Range("A266"))
$ws.Range("A266").Value = 44340
$ws.Range("B266").Value = 0
$ws.Range("C266").Value = 7
$ws.Range("D266").Value = 38.93647791745467

$ws.Range("A255").Copy($ws.Range("A267"))
$ws.Range("A267").Value = 44341
$ws.Range("B267").Value = 2
$ws.Range("C267").Value = 8
$ws.Range("D267").Value = 44.49883190566248

$ws.Range("A255").Copy($ws.Range("A268"))
$ws.Range("A268").Value = 44342
$ws.Range("B268").Value = 0
$ws.Range("C268").Value = 7
$ws.Range("D268").Value = 38.93647791745467

$ws.Range("A255").Copy($ws.Range("A269"))
$ws.Range("A269").Value = 44343
$ws.Range("B269").Value = 1
$ws.Range("C269").Value = 4
$ws.Range("D269").Value = 22.24941595283124
